$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 16, pushing existing rows 16-19 down to 17-20
$ws.Rows("16").Insert()

# Populate the newly inserted row 16 with the new record
$ws.Range("A16").Value = 1
$ws.Range("B16").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C16").Value = "Arica y Parinacota"
$ws.Range("D16").Value = 44435
$ws.Range("E16").Value = 15
$ws.Range("F16").Value = 100112044
$ws.Range("G16").Value = "Perejil"
$ws.Range("H16").Value = "Sin especificar"
$ws.Range("I16").Value = "Primera"
$ws.Range("J16").Value = 300
$ws.Range("K16").Value = 900
$ws.Range("L16").Value = 1000
$ws.Range("M16").Value = 950
$ws.Range("N16").Value = "$/atado 1,5 a 2 kilos"
$ws.Range("O16").Value = "Región de Arica y Parinacota"
$ws.Range("P16").Value = 475
$ws.Range("Q16").Value = 2
$ws.Range("R16").Value = "Hortaliza"

# Keep D16 consistent with the date-formatted style used by the other date cells
$ws.Range("D16").NumberFormat = $ws.Range("D17").NumberFormat
